# -----------------------------------------------------------------------
# "Added changes for clarity"
#
# - Adds four new highlighted clarification bullets (ilvl 4 / numId 3,
#   plus one ilvl 1 / numId 4) under the existing "Idea Genie" scope list
#   and the "Team composition" list.
# - Removes the stray "_GoBack" bookmark left behind in the
#   "Requirement Analysts" bullet.
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

function Get-ParagraphByText($doc, [string]$text) {
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text.TrimEnd() -eq $text) {
            return $para
        }
    }
    return $null
}

function Add-BulletAfter($doc, [string]$anchorText, [string]$newText, [bool]$demote) {
    $anchor = Get-ParagraphByText $doc $anchorText
    $rng = $anchor.Range
    $rng.Collapse(0)          # wdCollapseEnd
    $rng.InsertParagraphAfter()
    $newPara = Get-ParagraphByText $doc $anchorText
    $newPara = $newPara.Next()
    $newPara.Range.Text = $newText
    $newPara.Range.HighlightColorIndex = 7   # wdYellow
    if ($demote) {
        $newPara.Range.ListFormat.ListLevelNumber = 4
    }
}

# ---------------------------------------------------------------------
# Scope-of-product clarifications (new sub-bullets, one level deeper
# than their parent bullet, highlighted yellow).
# ---------------------------------------------------------------------
Add-BulletAfter $d "Addressing new ideas which can be registered and converted into a Product" `
    "Completely a new idea/innovation" $true

Add-BulletAfter $d "Registering PIs/PoV/PoC from Projects/Hackathon/other events and is accessible for reference" `
    "Best practices being followed in a project/Reusable Components registrations" $true

Add-BulletAfter $d "Ideathon/Hackathon can be conducted with the help of the Portal and the inbuilt Framework built" `
    "Portal should support for conducting online events" $true

# ---------------------------------------------------------------------
# Team-composition clarification (new bullet at the same level as its
# siblings, highlighted yellow).
# ---------------------------------------------------------------------
Add-BulletAfter $d "Testers." "Operation Engineer" $false

# ---------------------------------------------------------------------
# Drop the leftover "_GoBack" bookmark around "Requirement Analysts".
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
